# ModelRuns_RTP2025.xlsx edit:
# Add two new 2023 runs (TM160_IPA_26 and TM160_IPA_27) with the new CDAP,
# inserted as rows 38-39 (pushing the existing rows 38-45 down to 40-47).
# Both new runs are identical except for the network version / description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 38 (2025 Plan row). Inserting
# immediately below row 37 copies row 37's number formatting/style
# (s=1 / s=10 for J / s=44 for M,N / s=6 for Q,R,S) onto the new rows,
# matching the style pattern used by the new data rows.
$ws.Rows("38:39").Insert()

# --- Row 38: 2023_TM160_IPA_26 ---
$ws.Cells.Item(38, 1).Value2 = 2023
$ws.Cells.Item(38, 2).Value2 = "2023_TM160_IPA_26"
$ws.Cells.Item(38, 3).Value2 = "RTP2025_IP"
$ws.Cells.Item(38, 4).Value2 = "Past year"
$ws.Cells.Item(38, 5).Value2 = "New CDAP, WFH at ~25%, , trn hes=120, 0, 45"
$ws.Cells.Item(38, 6).Value2 = "petrale"
$ws.Cells.Item(38, 7).Value2 = "n/a"
$ws.Cells.Item(38, 8).Value2 = "current"
$ws.Cells.Item(38, 9).Value2 = "BlueprintNetworks_v09\net_2023_Blueprint"
$ws.Cells.Item(38, 10).Value2 = "model2-a"
$ws.Cells.Item(38, 11).Value2 = "https://app.asana.com/0/1204085012544660/1205764227468992/f"
$ws.Cells.Item(38, 12).Value2 = 17.77
$ws.Cells.Item(38, 13).Value2 = "na"
$ws.Cells.Item(38, 14).Value2 = "na"
$ws.Cells.Item(38, 15).Value2 = 0.94
$ws.Cells.Item(38, 16).Value2 = 0.855
$ws.Cells.Item(38, 17).Value2 = 120
$ws.Cells.Item(38, 18).Value2 = 0
$ws.Cells.Item(38, 19).Value2 = 45

# --- Row 39: 2023_TM160_IPA_27 (same run, network v10) ---
$ws.Cells.Item(39, 1).Value2 = 2023
$ws.Cells.Item(39, 2).Value2 = "2023_TM160_IPA_27"
$ws.Cells.Item(39, 3).Value2 = "RTP2025_IP"
$ws.Cells.Item(39, 4).Value2 = "Past year"
# Network (col I) filled before the description (col E) to match the
# shared-string insertion order of the original edit.
$ws.Cells.Item(39, 9).Value2 = "BlueprintNetworks_v10\net_2023_Blueprint"
$ws.Cells.Item(39, 5).Value2 = "New CDAP, WFH at ~25%, trn hes=120, 0, 45, network v10"
$ws.Cells.Item(39, 6).Value2 = "petrale"
$ws.Cells.Item(39, 7).Value2 = "n/a"
$ws.Cells.Item(39, 8).Value2 = "current"
$ws.Cells.Item(39, 10).Value2 = "model2-a"
$ws.Cells.Item(39, 11).Value2 = "https://app.asana.com/0/1204085012544660/1205764227468992/f"
$ws.Cells.Item(39, 12).Value2 = 17.77
$ws.Cells.Item(39, 13).Value2 = "na"
$ws.Cells.Item(39, 14).Value2 = "na"
$ws.Cells.Item(39, 15).Value2 = 0.94
$ws.Cells.Item(39, 16).Value2 = 0.855
$ws.Cells.Item(39, 17).Value2 = 120
$ws.Cells.Item(39, 18).Value2 = 0
$ws.Cells.Item(39, 19).Value2 = 45

# Leave the selection on the newly-added description cell, as in the
# final saved state of the workbook.
$ws.Range("E39").Select()
